$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowValues {
    param($wsArg, $rowNum, $values)
    $arr = New-Object "object[,]" 1,20
    for ($i = 0; $i -lt 20; $i++) {
        $arr[0,$i] = $values[$i]
    }
    $rangeRef = "A" + $rowNum + ":T" + $rowNum
    $wsArg.Range($rangeRef).Value = $arr
}

Set-RowValues $ws 2 @("ECs", "Col4a1", "Cd93", "ECs", 3, 1, 179.1580256666667, 537.4740770000001, 0.3468013736386751, 0.3468013736386751, 3, 1, 124.722578, 374.167734, 0.9767311432246923, 0.9767311432246923, 22345.05083053684, 201105.4574748315, 0.3387317021459968, 0.3387317021459968)
Set-RowValues $ws 3 @("ECs", "Col4a1", "Cd93", "FAPs", 3, 1, 179.1580256666667, 537.4740770000001, 0.3468013736386751, 0.3468013736386751, 1, 0.3333333333333333, 0.134712, 0.404136, 0.001054960600366076, 0.001054960600366076, 24.134735953608, 217.212623582472, 0.0003658617853416366, 0.0003658617853416365)
Set-RowValues $ws 4 @("ECs", "Col4a1", "Cd93", "sCs", 3, 1, 179.1580256666667, 537.4740770000001, 0.3468013736386751, 0.3468013736386751, 3, 1, 2.836578333333333, 8.509734999999999, 0.02221389617494163, 0.02221389617494163, 508.1957738488439, 4573.761964639595, 0.007703809707336667, 0.007703809707336666)
Set-RowValues $ws 5 @("FAPs", "Col4a1", "Cd93", "ECs", 3, 1, 239.807332, 719.421996, 0.4642019905988459, 0.4642019905988459, 3, 1, 124.722578, 374.167734, 0.9767311432246923, 0.9767311432246923, 29909.3886703419, 269184.498033077, 0.4534005409647887, 0.4534005409647887)
Set-RowValues $ws 6 @("FAPs", "Col4a1", "Cd93", "FAPs", 3, 1, 239.807332, 719.421996, 0.4642019905988459, 0.4642019905988459, 1, 0.3333333333333333, 0.134712, 0.404136, 0.001054960600366076, 0.001054960600366076, 32.304925308384, 290.744327775456, 0.0004897148106932861, 0.0004897148106932861)
Set-RowValues $ws 7 @("FAPs", "Col4a1", "Cd93", "sCs", 3, 1, 239.807332, 719.421996, 0.4642019905988459, 0.4642019905988459, 3, 1, 2.836578333333333, 8.509734999999999, 0.02221389617494163, 0.02221389617494163, 680.2322821256734, 6122.09053913106, 0.01031173482336399, 0.01031173482336399)
Set-RowValues $ws 8 @("sCs", "Col4a1", "Cd93", "ECs", 3, 1, 97.63589966666666, 292.907699, 0.1889966357624789, 0.1889966357624789, 3, 1, 124.722578, 374.167734, 0.9767311432246923, 0.9767311432246923, 12177.40111177601, 109596.6100059841, 0.1845989001139068, 0.1845989001139068)
Set-RowValues $ws 9 @("sCs", "Col4a1", "Cd93", "FAPs", 3, 1, 97.63589966666666, 292.907699, 0.1889966357624789, 0.1889966357624789, 1, 0.3333333333333333, 0.134712, 0.404136, 0.001054960600366076, 0.001054960600366076, 13.152727315896, 118.374545843064, 0.0001993840043311534, 0.0001993840043311534)
Set-RowValues $ws 10 @("sCs", "Col4a1", "Cd93", "sCs", 3, 1, 97.63589966666666, 292.907699, 0.1889966357624789, 0.1889966357624789, 3, 1, 2.836578333333333, 8.509734999999999, 0.02221389617494163, 0.02221389617494163, 276.9518775499739, 2492.566897949765, 0.004198351644240966, 0.004198351644240966)
